# Updated cryptos list (Price / Volume(1h) columns) per the GitHub Actions refresh.
# Note: some "Price" values look like plain decimals (e.g. 247.06) which Excel
# would otherwise auto-convert to a Number; a leading apostrophe forces them
# to stay Text, matching the original cell type/format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.032.46'
$ws.Range('E2').Value = '  -1.31%  '
$ws.Range('D3').Value = '2.171.47'
$ws.Range('E3').Value = '  -2.46%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '''247.06'
$ws.Range('E5').Value = '  -2.03%  '
$ws.Range('D6').Value = '''0.607'
$ws.Range('E6').Value = '  -3.45%  '
$ws.Range('D7').Value = '''66.38'
$ws.Range('E7').Value = '  -6.74%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '''0.565'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('D10').Value = '''59.65'
$ws.Range('E10').Value = '  +1.52%  '
$ws.Range('D11').Value = '''0.0929'
$ws.Range('E11').Value = '  -3.59%  '
$ws.Range('D12').Value = '''35.64'
$ws.Range('E12').Value = '  -16.35%  '
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('D14').Value = '''6.87'
$ws.Range('E14').Value = '  -2.66%  '
$ws.Range('D15').Value = '2.485.08'
$ws.Range('E15').Value = '  -2.95%  '
$ws.Range('D16').Value = '''0.855'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = '''14.24'
$ws.Range('E17').Value = '  -4.51%  '
$ws.Range('D18').Value = '2.170.13'
$ws.Range('E18').Value = '  -2.39%  '
$ws.Range('D19').Value = '40.931.27'
$ws.Range('E19').Value = '  -1.56%  '
$ws.Range('D20').Value = '0.0₃0939'
$ws.Range('E20').Value = '  -3.06%  '
$ws.Range('D21').Value = '''6.08'
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('D22').Value = '''71.37'
$ws.Range('E22').Value = '  -2.19%  '
$ws.Range('D23').Value = '''229.37'
$ws.Range('E23').Value = '  -2.21%  '
$ws.Range('D24').Value = '''2.09'
$ws.Range('E24').Value = '  -8.31%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = '''11.31'
$ws.Range('E26').Value = '  +8.44%  '
$ws.Range('D27').Value = '''3.67'
$ws.Range('E27').Value = '  -5.50%  '
$ws.Range('D28').Value = '''2.42'
$ws.Range('E28').Value = '  -3.41%  '
$ws.Range('E29').Value = '  -5.69%  '
$ws.Range('D30').Value = '''168.89'
$ws.Range('E30').Value = '  -1.60%  '
$ws.Range('D31').Value = '''1.99'
$ws.Range('E31').Value = '  -9.02%  '
$ws.Range('D32').Value = '''20.18'
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('E33').Value = '  -0.73%  '
$ws.Range('E34').Value = '  +1.72%  '
$ws.Range('D35').Value = '''0.0751'
$ws.Range('E35').Value = '  +4.11%  '
$ws.Range('E36').Value = '  -3.05%  '
$ws.Range('D37').Value = '''4.56'
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('D38').Value = '''3.96'
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('E39').Value = '  -8.04%  '
$ws.Range('E40').Value = '  +4.27%  '
$ws.Range('E41').Value = '  -4.91%  '
$ws.Range('D42').Value = '''5.46'
$ws.Range('E42').Value = '  -8.94%  '
$ws.Range('D43').Value = '''11.37'
$ws.Range('E43').Value = '  -7.83%  '
$ws.Range('D44').Value = '''4.86'
$ws.Range('E44').Value = '  -4.10%  '
$ws.Range('D45').Value = '''60.33'
$ws.Range('E45').Value = '  -12.85%  '
$ws.Range('E46').Value = '  -7.75%  '
$ws.Range('E47').Value = '  -0.26%  '
$ws.Range('D48').Value = '''8.41'
$ws.Range('E48').Value = '  -4.60%  '
$ws.Range('D49').Value = '''0.0989'
$ws.Range('E49').Value = '  -2.80%  '
$ws.Range('D50').Value = '''1.15'
$ws.Range('E50').Value = '  -0.92%  '
$ws.Range('E51').Value = '  -4.00%  '
